# Auto-generated edit script: update crypto price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.342.06"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.250.05"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'247.52"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").Value = "'0.633"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "'76.56"
$ws.Range("E7").Value = "  +6.40%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.636"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("D10").Value = "'40.58"
$ws.Range("E10").Value = "  +5.91%  "
$ws.Range("D11").Value = "'0.0954"
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("D12").Value = "'7.27"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("E13").Value = "  -1.42%  "
$ws.Range("D14").Value = "2.586.63"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").Value = "'14.97"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").Value = "2.267.07"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "42.244.34"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").Value = "0.0₃0983"
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("D20").Value = "'6.18"
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("D21").Value = "'71.76"
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").Value = "'232.06"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("D23").Value = "'2.20"
$ws.Range("E23").Value = "  -4.57%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("E25").Value = "  -5.20%  "
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("E27").Value = "  -4.27%  "
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("D29").Value = "'7.01"
$ws.Range("E29").Value = "  +8.95%  "
$ws.Range("D30").Value = "'168.13"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").Value = "'20.63"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").Value = "'0.0855"
$ws.Range("E32").Value = "  +6.65%  "
$ws.Range("D33").Value = "'31.97"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("D34").Value = "'0.121"
$ws.Range("E34").Value = "  -5.37%  "
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").Value = "'4.50"
$ws.Range("E36").Value = "  -6.20%  "
$ws.Range("D37").Value = "'4.79"
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("D38").Value = "'0.0299"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").Value = "'13.06"
$ws.Range("E39").Value = "  -6.75%  "
$ws.Range("E40").Value = "  -4.70%  "
$ws.Range("D41").Value = "'5.98"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "'117.14"
$ws.Range("E42").Value = "  +21.89%  "
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("D44").Value = "'60.35"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("E45").Value = "  -5.43%  "
$ws.Range("E46").Value = "  -1.78%  "
$ws.Range("D47").Value = "'0.997"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("E48").Value = "  -3.18%  "
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").Value = "'4.27"
$ws.Range("E50").Value = "  -12.56%  "
$ws.Range("D51").Value = "'4.16"
$ws.Range("E51").Value = "  -1.31%  "
